$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data - same as row 92 but with the next sequential date
$ws.Range("A93").Value = 45649
$ws.Range("B93").Value = 116.4121952
$ws.Range("C93").Value = 0.00170247
$ws.Range("D93").Value = 0.008850780000000001
$ws.Range("E93").Value = 0.06933635
$ws.Range("F93").Value = 12792.90181321
$ws.Range("G93").Value = 465.80531254
$ws.Range("H93").Value = 0.24
$ws.Range("I93").Value = 1.7904431
$ws.Range("J93").Value = 485.38834923

# Copy style/format of row 92 into row 93 last, so it is the one that sticks
$ws.Range("A92:J92").Copy() | Out-Null
$ws.Range("A93:J93").PasteSpecial(-4122) | Out-Null

